$d = $word.ActiveDocument

# Update the header date line.
$d.Content.Find.Execute("2025-09-25 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-09-26 Friday", 2) | Out-Null

# Update the multiplication answers in the table, cell-by-cell so that
# identical old/new text values in different cells never collide.
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Old = "48×40=1920"; New = "76×91=6916" },
    @{ Row = 1;  Col = 2; Old = "28×87=2436"; New = "37×96=3552" },
    @{ Row = 1;  Col = 3; Old = "65×78=5070"; New = "65×89=5785" },
    @{ Row = 1;  Col = 4; Old = "46×32=1472"; New = "92×82=7544" },
    @{ Row = 1;  Col = 5; Old = "73×91=6643"; New = "84×44=3696" },

    @{ Row = 5;  Col = 1; Old = "44×91=4004"; New = "44×99=4356" },
    @{ Row = 5;  Col = 2; Old = "35×91=3185"; New = "86×89=7654" },
    @{ Row = 5;  Col = 3; Old = "54×47=2538"; New = "46×19=874"  },
    @{ Row = 5;  Col = 4; Old = "68×20=1360"; New = "57×77=4389" },
    @{ Row = 5;  Col = 5; Old = "87×80=6960"; New = "54×22=1188" },

    @{ Row = 10; Col = 1; Old = "50×24=1200"; New = "16×57=912"  },
    @{ Row = 10; Col = 2; Old = "44×11=484";  New = "23×75=1725" },
    @{ Row = 10; Col = 3; Old = "46×46=2116"; New = "74×18=1332" },
    @{ Row = 10; Col = 4; Old = "35×40=1400"; New = "32×99=3168" },
    @{ Row = 10; Col = 5; Old = "90×27=2430"; New = "45×35=1575" },

    @{ Row = 15; Col = 1; Old = "87×34=2958"; New = "36×65=2340" },
    @{ Row = 15; Col = 2; Old = "86×66=5676"; New = "58×75=4350" },
    @{ Row = 15; Col = 3; Old = "83×24=1992"; New = "59×20=1180" },
    @{ Row = 15; Col = 4; Old = "13×72=936";  New = "28×87=2436" },
    @{ Row = 15; Col = 5; Old = "61×21=1281"; New = "99×40=3960" },

    @{ Row = 20; Col = 1; Old = "89×17=1513"; New = "48×11=528"  },
    @{ Row = 20; Col = 2; Old = "13×55=715";  New = "28×93=2604" },
    @{ Row = 20; Col = 3; Old = "85×32=2720"; New = "85×61=5185" },
    @{ Row = 20; Col = 4; Old = "44×41=1804"; New = "80×85=6800" },
    @{ Row = 20; Col = 5; Old = "71×21=1491"; New = "16×60=960"  }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    $rng.Find.Execute($u.Old, $true, $false, $false, $false, $false,
                       $true, 1, $false, $u.New, 2) | Out-Null
}
